$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Parqueados Primera Vez"  (A1:C5 -> A1:C6)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Update existing rows 3-5 with new plate/date data
$ws1.Range("A3").Value = 7
$ws1.Range("B3").Value = "ULM345"
$ws1.Range("C3").Value = "2024-02-27 12:43:53"

$ws1.Range("A4").Value = 11
$ws1.Range("B4").Value = "DML349"
$ws1.Range("C4").Value = "2024-02-27 15:12:38"

$ws1.Range("A5").Value = 13
$ws1.Range("B5").Value = "DMM974"
$ws1.Range("C5").Value = "2024-02-27 15:51:41"

# Add new row 6, cloning the formatting of row 5 first
$ws1.Range("A5:C5").Copy($ws1.Range("A6:C6"))
$ws1.Range("A6").Value = 14
$ws1.Range("B6").Value = "DMR978"
$ws1.Range("C6").Value = "2024-02-27 15:51:53"

# ---------------------------------------------------------------------------
# Sheet 2: "Vehiculos mas registrados (D.P)"  (A1:C9 -> A1:C12)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = 30
$ws2.Range("B3").Value = "MIK782"
$ws2.Range("C3").Value = 8

$ws2.Range("A4").Value = 29
$ws2.Range("B4").Value = "KLS878"
$ws2.Range("C4").Value = 7

$ws2.Range("A5").Value = 8
$ws2.Range("B5").Value = "ULM349"
$ws2.Range("C5").Value = 5

$ws2.Range("A6").Value = 6
$ws2.Range("B6").Value = "UML777"
$ws2.Range("C6").Value = 5

$ws2.Range("A7").Value = 4
$ws2.Range("B7").Value = "XML343"
$ws2.Range("C7").Value = 3

$ws2.Range("A8").Value = 10
$ws2.Range("B8").Value = "DML343"
$ws2.Range("C8").Value = 3

$ws2.Range("A9").Value = 12
$ws2.Range("B9").Value = "DML777"
$ws2.Range("C9").Value = 3

# Add rows 10-12, cloning the formatting of row 9 first
$ws2.Range("A9:C9").Copy($ws2.Range("A10:C10"))
$ws2.Range("A10").Value = 5
$ws2.Range("B10").Value = "ZAP234"
$ws2.Range("C10").Value = 2

$ws2.Range("A9:C9").Copy($ws2.Range("A11:C11"))
$ws2.Range("A11").Value = 13
$ws2.Range("B11").Value = "DMM974"
$ws2.Range("C11").Value = 1

$ws2.Range("A9:C9").Copy($ws2.Range("A12:C12"))
$ws2.Range("A12").Value = 14
$ws2.Range("B12").Value = "DMR978"
$ws2.Range("C12").Value = 1

# ---------------------------------------------------------------------------
# Sheet 3: "Vehiculos mas registrados en P"  (A1:C9 -> A1:C12)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "Vehículos más veces registrados en un parqueadero: Parqueadero La Tertulia"

$ws3.Range("A3").Value = 30
$ws3.Range("B3").Value = "MIK782"
$ws3.Range("C3").Value = 8

$ws3.Range("A4").Value = 6
$ws3.Range("B4").Value = "UML777"
$ws3.Range("C4").Value = 5

$ws3.Range("A5").Value = 8
$ws3.Range("B5").Value = "ULM349"
$ws3.Range("C5").Value = 5

$ws3.Range("A6").Value = 12
$ws3.Range("B6").Value = "DML777"
$ws3.Range("C6").Value = 3

$ws3.Range("A7").Value = 10
$ws3.Range("B7").Value = "DML343"
$ws3.Range("C7").Value = 2

$ws3.Range("A8").Value = 11
$ws3.Range("B8").Value = "DML349"
$ws3.Range("C8").Value = 1

$ws3.Range("A9").Value = 7
$ws3.Range("B9").Value = "ULM345"
$ws3.Range("C9").Value = 1

# Add rows 10-12, cloning the formatting of row 9 first
$ws3.Range("A9:C9").Copy($ws3.Range("A10:C10"))
$ws3.Range("A10").Value = 14
$ws3.Range("B10").Value = "DMR978"
$ws3.Range("C10").Value = 1

$ws3.Range("A9:C9").Copy($ws3.Range("A11:C11"))
$ws3.Range("A11").Value = 13
$ws3.Range("B11").Value = "DMM974"
$ws3.Range("C11").Value = 1

$ws3.Range("A9:C9").Copy($ws3.Range("A12:C12"))
$ws3.Range("A12").Value = 4
$ws3.Range("B12").Value = "XML343"
$ws3.Range("C12").Value = 1

# ---------------------------------------------------------------------------
# Sheet 4: "Ganancias de un parqueadero"  (A1:A5, unchanged extent)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("A1").Value = "Ganancias de un parqueadero: Parqueadero La Tertulia"
$ws4.Range("A2").Value = "Las ganancias de la fecha de hoy 2024-03-18 son: `$ 0"
$ws4.Range("A3").Value = "Las ganancias de esta semana son: `$ 0"
$ws4.Range("A4").Value = "Las ganancias del mes de MARCH son: `$ 2.400"
$ws4.Range("A5").Value = "Las ganancias del año de 2024 son: `$ 32.746"

# ---------------------------------------------------------------------------
# Sheet 5: "Coincidencias de placa"  (A1:C5, unchanged extent)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Coincidencias de la placa: DM"

$ws5.Range("A3").Value = 33
$ws5.Range("B3").Value = "DML349"
$ws5.Range("C3").Value = "2024-02-27 15:12:38"

$ws5.Range("A4").Value = 35
$ws5.Range("B4").Value = "DMM974"
$ws5.Range("C4").Value = "2024-02-27 15:51:41"

$ws5.Range("A5").Value = 36
$ws5.Range("B5").Value = "DMR978"
$ws5.Range("C5").Value = "2024-02-27 15:51:53"
